$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$cols = @("B","C","D","E","F","I","J","K","L","M","N")

$data = @{
  2 = @(1.02, 1.027703811178301, 1.032237439393364, 1.031335604262437, 1.03865560699846, 1.036071163284896, 1.032860207725328, 1.035043249014786, 1.03414401866725, 1.04144302081561, 1.014968057309961)
  3 = @(1.02, 1.028543015466689, 1.032859274340459, 1.032122164347624, 1.039763551559392, 1.036299139637528, 1.033340275995078, 1.035474314750352, 1.034739180058958, 1.042360236894961, 1.015128279350177)
  4 = @(1.02, 1.029086437180956, 1.033261910430599, 1.032631864755934, 1.040481195187561, 1.036445505994942, 1.033650664241427, 1.035752819254406, 1.035124377069898, 1.042953879975131, 1.015231834391668)
  5 = @(1.02, 1.02931498625713, 1.033431241166733, 1.032846319265842, 1.040783066570264, 1.036506762573411, 1.033781091110396, 1.035869800085334, 1.035286333720837, 1.043203481079524, 1.015275339993582)
  6 = @(1.02, 1.029353366162333, 1.03345967616612, 1.032882337425028, 1.040833762289327, 1.036517031629061, 1.033802986788405, 1.035889435644294, 1.03531352808326, 1.043245392151543, 1.01528264306793)
  7 = @(1.02, 1.029089490694823, 1.033264172792353, 1.03263472961614, 1.040485228124809, 1.036446325592974, 1.033652407249984, 1.035754382761845, 1.035126541066177, 1.042957215026288, 1.015232415829817)
  8 = @(1.02, 1.027987340519283, 1.032447535197405, 1.031601271077328, 1.039029890645705, 1.036148446614264, 1.033022499815649, 1.035189016863477, 1.034345137352378, 1.041752968283604, 1.015022229585292)
  9 = @(1.02, 1.026048334160984, 1.031010627347087, 1.029785946121627, 1.036471013803632, 1.035614768918399, 1.031910665741657, 1.034189567467052, 1.032968929089822, 1.039632055838935, 1.014650956905646)
  10 = @(1.02, 1.024757838129263, 1.030054195683375, 1.028579693005915, 1.034768905083044, 1.035253119542728, 1.031168252034112, 1.033521173297377, 1.032052014656906, 1.038218907988639, 1.014402860052835)
  11 = @(1.02, 1.024199570524284, 1.029640425871557, 1.028058331052813, 1.034032783375642, 1.035095138081956, 1.030846508150749, 1.033231266654074, 1.031655127557225, 1.037607195083734, 1.014295297947206)
  12 = @(1.02, 1.023992285280202, 1.029486790601156, 1.027864818787718, 1.03375949106036, 1.035036249252397, 1.03072695794221, 1.033123510130574, 1.031507728649891, 1.037380006662122, 1.014255324739321)
  13 = @(1.02, 1.024036745019786, 1.029519743273594, 1.027906321249748, 1.033818107025176, 1.0350488904858, 1.030752603672047, 1.033146627526404, 1.031539345173009, 1.0374287380295, 1.014263900020717)
  14 = @(1.02, 1.024182434618187, 1.029627725155639, 1.028042332306463, 1.034010190169302, 1.035090274542098, 1.030836626907953, 1.03322236093532, 1.031642943049131, 1.037588415023524, 1.014291994152285)
  15 = @(1.02, 1.024272209452166, 1.029694263988269, 1.028126152437044, 1.034128556977351, 1.035115745149289, 1.030888391097271, 1.033269013237962, 1.031706776153649, 1.037686801146559, 1.014309301254035)
  16 = @(1.02, 1.024794899691581, 1.030081664173001, 1.028614314301934, 1.034817778112042, 1.035263575123646, 1.031189599426407, 1.033540403256357, 1.03207835785669, 1.038259509395616, 1.014409995791786)
  17 = @(1.02, 1.025122911074216, 1.030324770543041, 1.028920781716111, 1.035250350293605, 1.035355934597059, 1.031378466798608, 1.033710509191193, 1.032311480322161, 1.038618805238114, 1.014473123050842)
  18 = @(1.02, 1.025314285203863, 1.030466606113558, 1.029099630881797, 1.035502749385477, 1.035409672653824, 1.031488603456088, 1.033809681998358, 1.03244747042197, 1.038828394666816, 1.014509931108651)
  19 = @(1.02, 1.025379547388748, 1.03051497439503, 1.029160629368672, 1.035588825705774, 1.035427973253419, 1.031526152689189, 1.033843489343368, 1.032493841821718, 1.038899862324042, 1.01452247948846)
  20 = @(1.02, 1.025087713310385, 1.030298683823174, 1.02888789114556, 1.035203930384265, 1.035346039120651, 1.031358205849708, 1.033692263299756, 1.032286467061377, 1.038580254303774, 1.01446635143588)
  21 = @(1.02, 1.024139530438904, 1.02959592559036, 1.028002276431552, 1.033953622711205, 1.035078093691781, 1.030811885261513, 1.033200061321783, 1.031612435422601, 1.037541393324168, 1.014283721675914)
  22 = @(1.02, 1.023543834860818, 1.029154405611196, 1.027446293697608, 1.033168292100123, 1.034908425635538, 1.03046815967723, 1.032890176642011, 1.031188776501673, 1.036888388173025, 1.014168780413467)
  23 = @(1.02, 1.023859579907832, 1.029388431747957, 1.027740950691372, 1.03358453592083, 1.034998483445524, 1.030650396804585, 1.033054491651826, 1.031413353256282, 1.037234542459233, 1.014229723702327)
  24 = @(1.02, 1.025103617508438, 1.030310471178739, 1.028902752694825, 1.035224905274132, 1.035350510873537, 1.031367360982324, 1.033700507978363, 1.032297769435005, 1.038597673755203, 1.01446941127742)
  25 = @(1.02, 1.02654923595189, 1.031381843222511, 1.030254559132037, 1.037131875979425, 1.035753773709582, 1.032198315618573, 1.034448322561398, 1.033324618755594, 1.04018022542988, 1.014747043922057)
}

foreach ($row in $data.Keys) {
    $rowVals = $data[$row]
    for ($i = 0; $i -lt $cols.Length; $i++) {
        $ws.Range($cols[$i] + $row).Value = $rowVals[$i]
    }
}